$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.153.10"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "2.048.74"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.36"
$ws.Range("E5").Value = "  -1.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.665"
$ws.Range("E6").Value = "  -0.92%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.00"
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.384"
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0784"
$ws.Range("E10").Value = "  -1.92%  "
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.22"
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.900"
$ws.Range("E13").Value = "  +11.55%  "
$ws.Range("D14").Value = "2.346.51"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.74"
$ws.Range("E15").Value = "  +1.75%  "
$ws.Range("D16").Value = "2.049.27"
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.97"
$ws.Range("E17").Value = "  +15.46%  "
$ws.Range("D18").Value = "37.157.17"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.58"
$ws.Range("E19").Value = "  -1.22%  "
$ws.Range("D20").Value = "0.0₃0897"
$ws.Range("E20").Value = "  -2.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.44"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "236.85"
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("E24").Value = "  +3.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.61"
$ws.Range("E25").Value = "  +3.03%  "
$ws.Range("E26").Value = "  -4.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.91"
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.18"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("E29").Value = "  -0.76%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.99"
$ws.Range("E30").Value = "  +4.70%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.17"
$ws.Range("E31").Value = "  +0.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0622"
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.52"
$ws.Range("E33").Value = "  +0.89%  "
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.26"
$ws.Range("E36").Value = "  -1.06%  "
$ws.Range("E37").Value = "  +1.18%  "
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.30"
$ws.Range("E39").Value = "  +14.21%  "
$ws.Range("E40").Value = "  +8.36%  "
$ws.Range("E41").Value = "  -14.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0224"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.60"
$ws.Range("E43").Value = "  -1.84%  "
$ws.Range("E44").Value = "  -1.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "96.18"
$ws.Range("E45").Value = "  -2.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.45"
$ws.Range("E46").Value = "  -1.41%  "
$ws.Range("D47").Value = "1.273.74"
$ws.Range("E47").Value = "  -1.66%  "
$ws.Range("E48").Value = "  -2.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.82"
$ws.Range("E49").Value = "  -1.48%  "
$ws.Range("D50").Value = "2.231.28"
$ws.Range("E50").Value = "  -1.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.46"
$ws.Range("E51").Value = "  +0.75%  "
